$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 42.409254
$ws.Range("H2").Value = 127.227762
$ws.Range("I2").Value = 0.6138221220752584
$ws.Range("J2").Value = 0.6138221220752584
$ws.Range("O2").Value = 0.8416031693647025
$ws.Range("P2").Value = 0.8416031693647025
$ws.Range("Q2").Value = 66.98032758251999
$ws.Range("R2").Value = 602.82294824268
$ws.Range("S2").Value = 0.5165946433647048
$ws.Range("T2").Value = 0.5165946433647048
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 42.409254
$ws.Range("H3").Value = 127.227762
$ws.Range("I3").Value = 0.6138221220752584
$ws.Range("J3").Value = 0.6138221220752584
$ws.Range("M3").Value = 0.2972526666666667
$ws.Range("N3").Value = 0.8917580000000001
$ws.Range("O3").Value = 0.1583968306352975
$ws.Range("P3").Value = 0.1583968306352975
$ws.Range("Q3").Value = 12.606263842844
$ws.Range("R3").Value = 113.456374585596
$ws.Range("S3").Value = 0.09722747871055362
$ws.Range("T3").Value = 0.09722747871055362
$ws.Range("I4").Value = 0.07014398987036251
$ws.Range("J4").Value = 0.07014398987036251
$ws.Range("O4").Value = 0.8416031693647025
$ws.Range("P4").Value = 0.8416031693647025
$ws.Range("S4").Value = 0.05903340418678267
$ws.Range("T4").Value = 0.05903340418678267
$ws.Range("I5").Value = 0.07014398987036251
$ws.Range("J5").Value = 0.07014398987036251
$ws.Range("M5").Value = 0.2972526666666667
$ws.Range("N5").Value = 0.8917580000000001
$ws.Range("O5").Value = 0.1583968306352975
$ws.Range("P5").Value = 0.1583968306352975
$ws.Range("Q5").Value = 1.440569851581778
$ws.Range("R5").Value = 12.965128664236
$ws.Range("S5").Value = 0.01111058568357983
$ws.Range("T5").Value = 0.01111058568357983
$ws.Range("G6").Value = 21.83492733333334
$ws.Range("H6").Value = 65.50478200000001
$ws.Range("I6").Value = 0.3160338880543792
$ws.Range("J6").Value = 0.3160338880543791
$ws.Range("O6").Value = 0.8416031693647025
$ws.Range("P6").Value = 0.8416031693647025
$ws.Range("Q6").Value = 34.48564753172
$ws.Range("R6").Value = 310.37082778548
$ws.Range("S6").Value = 0.2659751218132151
$ws.Range("T6").Value = 0.265975121813215
$ws.Range("G7").Value = 21.83492733333334
$ws.Range("H7").Value = 65.50478200000001
$ws.Range("I7").Value = 0.3160338880543792
$ws.Range("J7").Value = 0.3160338880543791
$ws.Range("M7").Value = 0.2972526666666667
$ws.Range("N7").Value = 0.8917580000000001
$ws.Range("O7").Value = 0.1583968306352975
$ws.Range("P7").Value = 0.1583968306352975
$ws.Range("Q7").Value = 6.490490376306223
$ws.Range("R7").Value = 58.41441338675601
$ws.Range("S7").Value = 0.05005876624116407
$ws.Range("T7").Value = 0.05005876624116405
